$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("M9").Value = -11
$ws.Range("K9").Value = 180
$ws.Range("I9").Value = 180
$ws.Range("H9").Value = 192.5
$ws.Range("M33").Value = 24.75
$ws.Range("I33").Value = 204.25
$ws.Range("H33").Value = 246.42857
$ws.Range("K33").Value = 204.25
$ws.Range("J40").Value = 1497.5
$ws.Range("L40").Value = 1497.5
$ws.Range("H40").Value = 3364.5
$ws.Range("K40").Value = 3737.9
$ws.Range("M40").Value = -3562.9
$ws.Range("N40").Value = -1847.5
$ws.Range("I40").Value = 3737.9
$ws.Range("H95").Value = 78668.336
$ws.Range("N95").Value = -84160.336
$ws.Range("L95").Value = 78668.336
$ws.Range("J95").Value = 78668.336
$ws.Range("H99").Value = 211
$ws.Range("K99").Value = 633
$ws.Range("I99").Value = 211
$ws.Range("M99").Value = 865
$ws.Range("H100").Value = 2663.9
$ws.Range("M100").Value = -1510.6
$ws.Range("K100").Value = 2051.6
$ws.Range("I100").Value = 2051.6
$ws.Range("N101").ClearContents()
$ws.Range("J101").Value = 0
$ws.Range("K101").Value = 645.75
$ws.Range("I101").Value = 215.25
$ws.Range("M101").Value = 976.25
$ws.Range("L101").Value = 0
$ws.Range("H101").Value = 215.25
$ws.Range("H114").Value = 95000
$ws.Range("N114").Value = -103678
$ws.Range("L114").Value = 95000
$ws.Range("J114").Value = 95000
$ws.Range("I132").Value = 25499.666
$ws.Range("M132").Value = -73968.99800000001
$ws.Range("L132").Value = 22583.1432
$ws.Range("K132").Value = 76498.99800000001
$ws.Range("N132").Value = -27643.1432
$ws.Range("H132").Value = 17636.938
$ws.Range("J132").Value = 7527.7144
$ws.Range("K137").Value = 4167.6
$ws.Range("H137").Value = 6844.207
$ws.Range("I137").Value = 1389.2
$ws.Range("M137").Value = -1617.6

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("M2").Value = -472.8570999999999
$ws.Range("I2").Value = 585.8570999999999
$ws.Range("J2").Value = 2109.3333
$ws.Range("K2").Value = 585.8570999999999
$ws.Range("L2").Value = 2109.3333
$ws.Range("H2").Value = 1042.9
$ws.Range("N2").Value = -2335.3333
$ws.Range("M34").Value = -99729
$ws.Range("K34").Value = 100000
$ws.Range("H34").Value = 100000
$ws.Range("I34").Value = 100000
$ws.Range("J36").Value = 30000
$ws.Range("H36").Value = 6324.875
$ws.Range("L36").Value = 30000
$ws.Range("I36").Value = 2942.7144
$ws.Range("K36").Value = 2942.7144
$ws.Range("N36").Value = -30692
$ws.Range("M36").Value = -2596.7144
$ws.Range("H61").Value = 5681.5884
$ws.Range("M61").Value = -2140.923
$ws.Range("I61").Value = 2352.923
$ws.Range("K61").Value = 2352.923
$ws.Range("I74").Value = 897.89795
$ws.Range("H74").Value = 4271.2095
$ws.Range("M74").Value = -23.89795000000004
$ws.Range("K74").Value = 897.89795
$ws.Range("L74").Value = 16986
$ws.Range("N74").Value = -18734
$ws.Range("J74").Value = 16986
$ws.Range("J77").Value = 16986
$ws.Range("M77").Value = -121.4897500000006
$ws.Range("N77").Value = -93666
$ws.Range("K77").Value = 4489.489750000001
$ws.Range("I77").Value = 897.89795
$ws.Range("L77").Value = 84930
$ws.Range("H77").Value = 4271.2095
$ws.Range("N116").Value = -6697.3333
$ws.Range("K116").Value = 585.8570999999999
$ws.Range("H116").Value = 1042.9
$ws.Range("M116").Value = 1708.1429
$ws.Range("J116").Value = 2109.3333
$ws.Range("I116").Value = 585.8570999999999
$ws.Range("L116").Value = 2109.3333
$ws.Range("I132").Value = 5002499
$ws.Range("M132").Value = -15004967
$ws.Range("L132").Value = 0
$ws.Range("K132").Value = 15007497
$ws.Range("N132").ClearContents()
$ws.Range("H132").Value = 5002499
$ws.Range("J132").Value = 0
$ws.Range("I136").Value = 2352.923
$ws.Range("H136").Value = 5681.5884
$ws.Range("K136").Value = 7058.768999999999
$ws.Range("M136").Value = -4508.768999999999
$ws.Range("H141").Value = 0
$ws.Range("L141").Value = 0
$ws.Range("N141").ClearContents()
$ws.Range("J141").Value = 0

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1042.9
$ws.Range("J3").Value = 2109.3333
$ws.Range("I3").Value = 585.8570999999999
$ws.Range("M3").Value = -471.8570999999999
$ws.Range("K3").Value = 585.8570999999999
$ws.Range("N3").Value = -2337.3333
$ws.Range("L3").Value = 2109.3333
$ws.Range("I10").Value = 485
$ws.Range("H10").Value = 2994
$ws.Range("K10").Value = 485
$ws.Range("M10").Value = -345
$ws.Range("H44").Value = 2000
$ws.Range("I44").Value = 2000
$ws.Range("M44").Value = -1503
$ws.Range("K44").Value = 2000
$ws.Range("M94").Value = -1230.95
$ws.Range("H94").Value = 1825.0465
$ws.Range("K94").Value = 1681.95
$ws.Range("I94").Value = 1681.95
$ws.Range("L99").Value = 1999
$ws.Range("J99").Value = 1999
$ws.Range("H99").Value = 15926.714
$ws.Range("K99").Value = 18248
$ws.Range("I99").Value = 18248
$ws.Range("M99").Value = -16750
$ws.Range("N99").Value = -4995
$ws.Range("I107").Value = 769
$ws.Range("M107").Value = 1151
$ws.Range("K107").Value = 769
$ws.Range("H107").Value = 788.25
$ws.Range("K134").Value = 5245.950000000001
$ws.Range("I134").Value = 1748.65
$ws.Range("M134").Value = -2710.950000000001
$ws.Range("H134").Value = 6427.2856

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("K58").Value = 5186.75
$ws.Range("H58").Value = 13013.462
$ws.Range("M58").Value = -4983.75
$ws.Range("I58").Value = 5186.75
$ws.Range("I105").Value = 10665.2
$ws.Range("H105").Value = 7623.467
$ws.Range("K105").Value = 10665.2
$ws.Range("M105").Value = -8918.200000000001
$ws.Range("I136").Value = 5186.75
$ws.Range("H136").Value = 13013.462
$ws.Range("K136").Value = 15560.25
$ws.Range("M136").Value = -13010.25

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H32").Value = 2081.875
$ws.Range("N32").Value = -10745.9999
$ws.Range("I32").Value = 1295
$ws.Range("J32").Value = 3393.3333
$ws.Range("L32").Value = 10179.9999
$ws.Range("K32").Value = 3885
$ws.Range("M32").Value = -3602
$ws.Range("L129").Value = 15888.4995
$ws.Range("N129").Value = -25888.4995
$ws.Range("J129").Value = 5296.1665
$ws.Range("H129").Value = 1116292.8
$ws.Range("L131").Value = 18412.9659
$ws.Range("J131").Value = 6137.6553
$ws.Range("H131").Value = 6137.6553
$ws.Range("N131").Value = -28492.9659

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("M2").Value = -7692243
$ws.Range("I2").Value = 7692356
$ws.Range("K2").Value = 7692356
$ws.Range("H2").Value = 6250126.5
$ws.Range("J6").Value = 1500
$ws.Range("L6").Value = 1500
$ws.Range("H6").Value = 1250
$ws.Range("N6").Value = -1726
$ws.Range("J10").Value = 6971.143
$ws.Range("H10").Value = 11533.444
$ws.Range("L10").Value = 6971.143
$ws.Range("N10").Value = -7309.143
$ws.Range("K11").Value = 9460000
$ws.Range("I11").Value = 9460000
$ws.Range("H11").Value = 7966666.5
$ws.Range("M11").Value = -9459861
$ws.Range("J16").Value = 1500
$ws.Range("H16").Value = 1250
$ws.Range("L16").Value = 1500
$ws.Range("N16").Value = -2000
$ws.Range("K22").Value = 1500
$ws.Range("I22").Value = 1500
$ws.Range("H22").Value = 3500
$ws.Range("M22").Value = -971
$ws.Range("I102").Value = 2319.35
$ws.Range("K102").Value = 2319.35
$ws.Range("M102").Value = -697.3499999999999
$ws.Range("H102").Value = 2319.35
$ws.Range("N117").Value = -96184.25
$ws.Range("L117").Value = 89300.25
$ws.Range("H117").Value = 89300.25
$ws.Range("J117").Value = 89300.25

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("L22").Value = 1989.6666
$ws.Range("N22").Value = -2579.6666
$ws.Range("J22").Value = 1989.6666
$ws.Range("H22").Value = 1944.8387
$ws.Range("L27").Value = 1989.6666
$ws.Range("H27").Value = 1944.8387
$ws.Range("J27").Value = 1989.6666
$ws.Range("N27").Value = -2203.6666
$ws.Range("L46").Value = 7114.143
$ws.Range("J46").Value = 7114.143
$ws.Range("N46").Value = -7490.143
$ws.Range("H46").Value = 6337.375
$ws.Range("I55").Value = 2542.1428
$ws.Range("J55").Value = 1664.65
$ws.Range("N55").Value = -2010.65
$ws.Range("K55").Value = 2542.1428
$ws.Range("M55").Value = -2369.1428
$ws.Range("H55").Value = 1892.1482
$ws.Range("L55").Value = 1664.65
$ws.Range("I136").Value = 2401.5557
$ws.Range("H136").Value = 11364.909
$ws.Range("K136").Value = 7204.6671
$ws.Range("M136").Value = -4654.6671

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H3").Value = 837416.5
$ws.Range("I3").Value = 1252624.8
$ws.Range("M3").Value = -1252510.8
$ws.Range("K3").Value = 1252624.8
$ws.Range("J6").Value = 2000
$ws.Range("L6").Value = 2000
$ws.Range("H6").Value = 2000
$ws.Range("N6").Value = -2230
$ws.Range("K11").Value = 20000000
$ws.Range("N11").Value = -16673584
$ws.Range("J11").Value = 16673300
$ws.Range("I11").Value = 20000000
$ws.Range("L11").Value = 16673300
$ws.Range("H11").Value = 17504976
$ws.Range("M11").Value = -19999858
$ws.Range("I132").Value = 3913.8462
$ws.Range("M132").Value = -9211.5386
$ws.Range("L132").Value = 18000
$ws.Range("K132").Value = 11741.5386
$ws.Range("N132").Value = -23060
$ws.Range("H132").Value = 4062.8572
$ws.Range("J132").Value = 6000
$ws.Range("I136").Value = 4508.625
$ws.Range("H136").Value = 4529.077
$ws.Range("K136").Value = 13525.875
$ws.Range("M136").Value = -10975.875
